$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy style9 (fontId5+right) onto M7 scratch via paste
$ws.Range("C2").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = "gering"

# Step 2: Now I need to ALSO add yellow fill (style18 target) to M7 without breaking font.
# Try Interior.Color set directly (non-font non-alignment property) and see if it preserves existing style.
$ws.Range("M7").Interior.Color = 65535
Write-Host "Done"
